$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme
$cs.Colors(1).RGB = 0
Write-Host "done"
